$wb = $excel.ActiveWorkbook

# --- NextBus1 sheet updates (refreshed ETAs, statuses, and minutes-to-arrival) ---
$ws = $wb.Worksheets.Item("NextBus1")
$ws.Range("F2").Value = 45685.60255787037
$ws.Range("L2").Value = "SD"
$ws.Range("O2").Value = 14
$ws.Range("F3").Value = 45685.59535879629
$ws.Range("O3").Value = 4
$ws.Range("F4").Value = 45685.59855324074
$ws.Range("O4").Value = 8
$ws.Range("F5").Value = 45685.59276620371
$ws.Range("O5").Value = 0
$ws.Range("F6").Value = 45685.59403935185
$ws.Range("O6").Value = 2
$ws.Range("F7").Value = 45685.59466435185
$ws.Range("F8").Value = 45685.5975
$ws.Range("O8").Value = 7
$ws.Range("F9").Value = 45685.59313657408
$ws.Range("O9").Value = 0
$ws.Range("F10").Value = 45685.59788194444
$ws.Range("L10").Value = "DD"
$ws.Range("O10").Value = 7
$ws.Range("F11").Value = 45685.59325231481
$ws.Range("O11").Value = 1
$ws.Range("F12").Value = 45685.600625
$ws.Range("O12").Value = 11
$ws.Range("F13").Value = 45685.601875
$ws.Range("O13").Value = 13
$ws.Range("F14").Value = 45685.59815972222
$ws.Range("O14").Value = 8
$ws.Range("F15").Value = 45685.59827546297
$ws.Range("O15").Value = 8

# --- NextBus2 sheet updates ---
$ws = $wb.Worksheets.Item("NextBus2")
$ws.Range("F2").Value = 45685.60777777778
$ws.Range("L2").Value = "DD"
$ws.Range("O2").Value = 21
$ws.Range("F3").Value = 45685.59921296296
$ws.Range("L3").Value = "BD"
$ws.Range("O3").Value = 9
$ws.Range("F4").Value = 45685.6055324074
$ws.Range("O4").Value = 18
$ws.Range("F5").Value = 45685.59501157407
$ws.Range("O5").Value = 3
$ws.Range("F6").Value = 45685.60196759259
$ws.Range("O6").Value = 13
$ws.Range("F7").Value = 45685.60040509259
$ws.Range("O7").Value = 11
$ws.Range("F8").Value = 45685.60643518518
$ws.Range("O8").Value = 20
$ws.Range("F9").Value = 45685.59864583334
$ws.Range("O9").Value = 8
$ws.Range("F10").Value = 45685.60280092592
$ws.Range("L10").Value = "SD"
$ws.Range("O10").Value = 14
$ws.Range("F11").Value = 45685.59997685185
$ws.Range("O11").Value = 10
$ws.Range("F12").Value = 45685.61018518519
$ws.Range("J12").Value = 0
$ws.Range("O12").Value = 25
$ws.Range("F13").Value = 45685.61020833333
$ws.Range("L13").Value = "DD"
$ws.Range("O13").Value = 25
$ws.Range("F14").Value = 45685.60461805556
$ws.Range("O14").Value = 17
$ws.Range("F15").Value = 45685.60622685185
$ws.Range("O15").Value = 19

# --- NextBus3 sheet updates ---
$ws = $wb.Worksheets.Item("NextBus3")
$ws.Range("F2").Value = 45685.61422453704
$ws.Range("J2").Value = 0
$ws.Range("O2").Value = 31
$ws.Range("F3").Value = 45685.60741898148
$ws.Range("L3").Value = "SD"
$ws.Range("O3").Value = 21
$ws.Range("F4").Value = 45685.6131712963
$ws.Range("O4").Value = 29
$ws.Range("F5").Value = 45685.60120370371
$ws.Range("O5").Value = 12
$ws.Range("F6").Value = 45685.6055324074
$ws.Range("O6").Value = 18
$ws.Range("O7").Value = 19
$ws.Range("F8").Value = 45685.61542824074
$ws.Range("O8").Value = 32
$ws.Range("F9").Value = 45685.60469907407
$ws.Range("O9").Value = 17
$ws.Range("F10").Value = 45685.60837962963
$ws.Range("O10").Value = 22
$ws.Range("F11").Value = 45685.60335648148
$ws.Range("O11").Value = 15
$ws.Range("F12").Value = 45685.62255787037
$ws.Range("L12").Value = "DD"
$ws.Range("O12").Value = 43
$ws.Range("F13").Value = 45685.61799768519
$ws.Range("I13").Value = "SDA"
$ws.Range("L13").Value = "SD"
$ws.Range("O13").Value = 36
$ws.Range("F14").Value = 45685.61204861111
$ws.Range("O14").Value = 28
$ws.Range("F15").Value = 45685.6128587963
$ws.Range("O15").Value = 29
